$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 corresponds to ICSA_thou / Initial Jobless Claims - auto-updated data refresh
$ws.Range("E9").Value = 214000
$ws.Range("G9").Value = 364752.8735632184
$ws.Range("H9").Value = -8000
$ws.Range("I9").Value = -0.03603603603603604
